$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 370, shifting existing rows 370-438 down to 371-439
$ws.Rows.Item(370).Insert()

# Populate the newly inserted row 370 with the new data record
$ws.Range("A370").Value = 6
$ws.Range("B370").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C370").Value = "Metropolitana"
$ws.Range("D370").Value = 44474
$ws.Range("D370").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E370").Value = 13
$ws.Range("F370").Value = 100112021
$ws.Range("G370").Value = "Ají"
$ws.Range("H370").Value = "Inferno"
$ws.Range("I370").Value = "Segunda"
$ws.Range("J370").Value = 40
$ws.Range("K370").Value = 40000
$ws.Range("L370").Value = 45000
$ws.Range("M370").Value = 42875
$ws.Range("N370").Value = "$/caja 15 kilos"
$ws.Range("O370").Value = "Provincia de Huasco"
$ws.Range("P370").Value = 2858
$ws.Range("Q370").Value = 15
$ws.Range("R370").Value = "Hortaliza"
